$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estimation": a couple of hour estimates were corrected.
# ---------------------------------------------------------------------------
$wsEstimation = $wb.Worksheets.Item("Estimation")
$wsEstimation.Range("C29").Value = 5
$wsEstimation.Range("C30").Value = 5

# ---------------------------------------------------------------------------
# Sheet "Iteration #1": log of completed tasks was filled in for rows 14-20.
# New shared strings must be introduced in the exact order Excel would have
# written them (date first, then description, then duration) so the cells
# that reuse an already-seen value don't mint a duplicate.
# ---------------------------------------------------------------------------
$wsIter1 = $wb.Worksheets.Item("Iteration #1")

$d20170130 = Get-Date -Year 2017 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0
$d20170131 = Get-Date -Year 2017 -Month 1 -Day 31 -Hour 0 -Minute 0 -Second 0
$d20170206 = Get-Date -Year 2017 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0
$d20170207 = Get-Date -Year 2017 -Month 2 -Day 7 -Hour 0 -Minute 0 -Second 0
$d20170212 = Get-Date -Year 2017 -Month 2 -Day 12 -Hour 0 -Minute 0 -Second 0

# Row 14 already carries the date-formatted style; just fill the values.
$wsIter1.Range("A14").Value = $d20170130
$wsIter1.Range("B14").Value = "realisation de mon menu de slide"
$wsIter1.Range("C14").Value = "3h"

$wsIter1.Range("C15").Value = "2h"
$wsIter1.Range("B16").Value = "affichage des musique de mon telephone dans une activité"
$wsIter1.Range("B17").Value = "la musique demarre suit à un click"
$wsIter1.Range("B15").Value = "realisation de l'interface principale de mon application, test de fonctionement"
$wsIter1.Range("B19").Value = "Lecture sur le fonctionement d'un service "
$wsIter1.Range("B20").Value = "Test personnels sur le fonctionement d'un service et son utilisation"
$wsIter1.Range("B18").Value = "Controle de la lecture de la musique( l'utilisateur peut aller a la musique suivante, precedente,) "

$wsIter1.Range("C16").Value = "2h"
$wsIter1.Range("C17").Value = "2h"
$wsIter1.Range("C18").Value = "3h"
$wsIter1.Range("C19").Value = "2h"
$wsIter1.Range("C20").Value = "2h"

# Rows 15-17 adopt the same date style already used by A14 (numFmt 14, bordered,
# left aligned) - copy formats so no redundant style entry is created.
$wsIter1.Range("A14").Copy()
$wsIter1.Range("A15:A17").PasteSpecial(-4122)

$wsIter1.Range("A15").Value = $d20170130
$wsIter1.Range("A16").Value = $d20170131
$wsIter1.Range("A17").Value = $d20170131

# Rows 18-20 use a new date style (numFmt 14, bordered, no forced alignment).
$wsIter1.Range("A18").NumberFormat = "mm-dd-yy"
$wsIter1.Range("A18").Copy()
$wsIter1.Range("A19:A20").PasteSpecial(-4122)

$wsIter1.Range("A18").Value = $d20170206
$wsIter1.Range("A19").Value = $d20170207
$wsIter1.Range("A20").Value = $d20170212

$wsIter1.Range("B20").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Iteration #2": only the remembered selection changed.
# ---------------------------------------------------------------------------
$wsIter2 = $wb.Worksheets.Item("Iteration #2")
$wsIter2.Range("A14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore "Estimation" as the active sheet/selection, matching tabSelected.
# ---------------------------------------------------------------------------
$wsEstimation.Activate() | Out-Null
$wsEstimation.Range("B8").Select() | Out-Null
